$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header in column H, reusing the same header formatting as
# the existing header cells (bold, bordered, centered).
$ws.Cells.Item(1, 7).Copy()
$ws.Cells.Item(1, 8).PasteSpecial(-4122)
$ws.Cells.Item(1, 8).Value = "Save"

# Populate the Save column for every data row. Row 41 (2024-04-30) is
# flagged as saved (1); every other row is 0.
for ($r = 2; $r -le 52; $r++) {
    if ($r -eq 41) {
        $ws.Cells.Item($r, 8).Value = 1
    } else {
        $ws.Cells.Item($r, 8).Value = 0
    }
}
